$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '331.51'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-0.10%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '41.74'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '6.49%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.696'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-0.18%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.08350'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '4.04%'
$ws.Range("B6").Value = 'KuCoinToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '8.827'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '2.55%'
$ws.Range("B7").Value = 'FTXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '2.025'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '3.83%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '4.541'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '1.10%'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '2.32%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9299'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '0.95%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1291'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '3.74%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.1961'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '1.27%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09427'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '1.85%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.03916'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '9.85%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.1061'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '0.98%'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '1.00%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.006162'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-1.40%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.445'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '2.32%'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '2.32%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.254'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-4.91%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1371'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-0.25%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.2482'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-8.06%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04413'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-0.62%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001248'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-1.09%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004392'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-1.29%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001201'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-0.12%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02802'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '10.39%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05550'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '1.22%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007794'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '3.63%'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '2.48%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.008932'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-9.83%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002246'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '6.37%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.01113'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-6.84%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00007035'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '3.14%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000750'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.21%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.003503'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '14.22%'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-0.42%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002100'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.21%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002000'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.21%'
